$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CTAS")

# Row 4 - Inventory
$ws.Range("B4").Value = 533000000.0
$ws.Range("C4").Value = 534000000.0
$ws.Range("D4").Value = 488000000.0
$ws.Range("E4").Value = 409000000.0
$ws.Range("F4").Value = 353000000.0

# Row 13 - Accounts Payable
$ws.Range("B13").Value = 238000000.0
$ws.Range("C13").Value = 274000000.0
$ws.Range("D13").Value = 253000000.0
$ws.Range("E13").Value = 231000000.0
$ws.Range("F13").Value = 243000000.0

# Row 21 - Long Term Tax Liability (Deferred)
$ws.Range("B21").Value = 390000000.0
$ws.Range("C21").Value = 376000000.0
$ws.Range("D21").Value = 385000000.0
$ws.Range("E21").Value = 389000000.0
$ws.Range("F21").Value = 424000000.0

$wb.Save()
